$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G10").Value = 5.0
$ws.Range("G12").Value = 1.0
$ws.Range("G13").Value = 5.0
$ws.Range("G17").Value = 5.0
$ws.Range("G22").Value = 5.0
